$wb = $excel.ActiveWorkbook

# Sheet "Hoja1": update the conversion note text in A1
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.98 = 24364.72 pesos`n✅ 24364.72 pesos = 5.97 = 951.6 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Sheet "tasas": update the rate figures
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 167.25
$ws2.Range("O10").Value = 4075
$ws2.Range("N12").Value = 4080
$ws2.Range("O12").Value = 159.35
